# Auto-generated edit script: updates Leve profit calculation cells (columns H-N)
# across multiple sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 17016.143
$ws.Range("I70").Value = 3150
$ws.Range("J70").Value = 35504.332
$ws.Range("K70").Value = 9450
$ws.Range("L70").Value = 106512.996
$ws.Range("M70").Value = -9180
$ws.Range("N70").Value = -107052.996
# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 17016.143
$ws.Range("I73").Value = 3150
$ws.Range("J73").Value = 35504.332
$ws.Range("K73").Value = 9450
$ws.Range("L73").Value = 106512.996
$ws.Range("M73").Value = -8514
$ws.Range("N73").Value = -108384.996
# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 1956.1666
$ws.Range("I86").Value = 2107.4
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 2107.4
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -984.4000000000001
$ws.Range("N86").Value = -3446
# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 1956.1666
$ws.Range("I89").Value = 2107.4
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 10537
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = -4921
$ws.Range("N89").Value = -17232
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4308.375
$ws.Range("I138").Value = 2699.4
$ws.Range("J138").Value = 6990
$ws.Range("K138").Value = 8098.200000000001
$ws.Range("L138").Value = 20970
$ws.Range("M138").Value = -2958.200000000001
$ws.Range("N138").Value = -31250

$ws = $wb.Worksheets.Item("ARM")
# Row 8 (Leve Item ID 3011)
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 5999.6665
$ws.Range("I102").Value = 5999.6665
$ws.Range("K102").Value = 5999.6665
$ws.Range("M102").Value = -4377.6665
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045
# Row 133 (Leve Item ID 41857)
$ws.Range("H133").Value = 89930.664
$ws.Range("J133").Value = 99916.8
$ws.Range("L133").Value = 99916.8
$ws.Range("N133").Value = -104976.8

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 26814
$ws.Range("I86").Value = 2332.6667
$ws.Range("J86").Value = 41502.8
$ws.Range("K86").Value = 2332.6667
$ws.Range("L86").Value = 41502.8
$ws.Range("M86").Value = -1209.6667
$ws.Range("N86").Value = -43748.8
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 26814
$ws.Range("I89").Value = 2332.6667
$ws.Range("J89").Value = 41502.8
$ws.Range("K89").Value = 11663.3335
$ws.Range("L89").Value = 207514
$ws.Range("M89").Value = -6047.333500000001
$ws.Range("N89").Value = -218746
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 5224.75
$ws.Range("I105").Value = 6000
$ws.Range("K105").Value = 6000
$ws.Range("M105").Value = -4253
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3575.25
$ws.Range("I107").Value = 1433.6666
$ws.Range("K107").Value = 1433.6666
$ws.Range("M107").Value = 486.3334
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3318.625
$ws.Range("I134").Value = 3318.625
$ws.Range("K134").Value = 9955.875
$ws.Range("M134").Value = -7420.875
# Row 137 (Leve Item ID 42153)
$ws.Range("H137").Value = 99995
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("N137").Value = -110195

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2611.1052
$ws.Range("I31").Value = 848.44446
$ws.Range("J31").Value = 4197.5
$ws.Range("K31").Value = 848.44446
$ws.Range("L31").Value = 4197.5
$ws.Range("M31").Value = -553.44446
$ws.Range("N31").Value = -4787.5
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2611.1052
$ws.Range("I34").Value = 848.44446
$ws.Range("J34").Value = 4197.5
$ws.Range("K34").Value = 848.44446
$ws.Range("L34").Value = 4197.5
$ws.Range("M34").Value = -646.44446
$ws.Range("N34").Value = -4601.5

$ws = $wb.Worksheets.Item("CUL")
# Row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 72.454544
$ws.Range("I38").Value = 66
$ws.Range("J38").Value = 74.875
$ws.Range("K38").Value = 198
$ws.Range("L38").Value = 224.625
$ws.Range("M38").Value = 149
$ws.Range("N38").Value = -918.625
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 5199.9
$ws.Range("I80").Value = 1999
$ws.Range("K80").Value = 5997
$ws.Range("M80").Value = -5061
# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 5199.9
$ws.Range("I83").Value = 1999
$ws.Range("K83").Value = 17991
$ws.Range("M83").Value = -13311
# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 4000
$ws.Range("I129").Value = 5000
$ws.Range("K129").Value = 15000
$ws.Range("M129").Value = -10000
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2137.25
$ws.Range("I131").Value = 749.5
$ws.Range("J131").Value = 2599.8333
$ws.Range("K131").Value = 2248.5
$ws.Range("L131").Value = 7799.499899999999
$ws.Range("M131").Value = 2791.5
$ws.Range("N131").Value = -17879.4999

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 3374.5
$ws.Range("I70").Value = 3374.5
$ws.Range("K70").Value = 3374.5
$ws.Range("M70").Value = -3104.5
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 3374.5
$ws.Range("I73").Value = 3374.5
$ws.Range("K73").Value = 3374.5
$ws.Range("M73").Value = -2438.5
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3281.3333
$ws.Range("J80").Value = 3999
$ws.Range("L80").Value = 3999
$ws.Range("N80").Value = -5995
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3281.3333
$ws.Range("J83").Value = 3999
$ws.Range("L83").Value = 19995
$ws.Range("N83").Value = -29979
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2944.3572
$ws.Range("I102").Value = 3268.3635
$ws.Range("J102").Value = 1756.3334
$ws.Range("K102").Value = 3268.3635
$ws.Range("L102").Value = 1756.3334
$ws.Range("M102").Value = -1646.3635
$ws.Range("N102").Value = -5000.3334
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 1300.3334
$ws.Range("J107").Value = 3100
$ws.Range("L107").Value = 3100
$ws.Range("N107").Value = -6940

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 83335660
$ws.Range("I16").Value = 250000000
$ws.Range("K16").Value = 250000000
$ws.Range("M16").Value = -249999830
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2876.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2876.5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2876.5
$ws.Range("N22").Value = -3466.5
$ws.Range("M22").ClearContents()
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2876.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2876.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2876.5
$ws.Range("N27").Value = -3090.5
$ws.Range("M27").ClearContents()
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 4260.4
$ws.Range("I46").Value = 2522.75
$ws.Range("J46").Value = 4892.273
$ws.Range("K46").Value = 2522.75
$ws.Range("L46").Value = 4892.273
$ws.Range("M46").Value = -2334.75
$ws.Range("N46").Value = -5268.273
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 4000.8572
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 5142.857
$ws.Range("I68").Value = 2666.6667
$ws.Range("K68").Value = 2666.6667
$ws.Range("M68").Value = -1917.6667
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 5142.857
$ws.Range("I71").Value = 2666.6667
$ws.Range("K71").Value = 13333.3335
$ws.Range("M71").Value = -9589.333500000001
# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 2143.2666
$ws.Range("I82").Value = 1740.909
$ws.Range("J82").Value = 3249.75
$ws.Range("K82").Value = 1740.909
$ws.Range("L82").Value = 3249.75
$ws.Range("M82").Value = -1379.909
$ws.Range("N82").Value = -3971.75
# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 2143.2666
$ws.Range("I85").Value = 1740.909
$ws.Range("J85").Value = 3249.75
$ws.Range("K85").Value = 1740.909
$ws.Range("L85").Value = 3249.75
$ws.Range("M85").Value = -492.9090000000001
$ws.Range("N85").Value = -5745.75
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 15368
$ws.Range("I93").Value = 15368
$ws.Range("K93").Value = 15368
$ws.Range("M93").Value = -14120

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 3747.182
$ws.Range("I96").Value = 3027.375
$ws.Range("K96").Value = 3027.375
$ws.Range("M96").Value = -1654.375
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 724.86664
$ws.Range("I113").Value = 747.8333
$ws.Range("K113").Value = 2243.4999
$ws.Range("M113").Value = -73.4998999999998
